# Refresh the cryptos list: update each changed Price (D) / Volume(1h) (E)
# cell (and the three re-ranked Coin/Link/Price/Volume rows) to match the
# latest coinranking.com snapshot pulled by the GitHub Actions job.
#
# All source values are plain text in the sheet (inline strings), including
# ones that look like plain decimals (e.g. "226.49"). Excel's COM Range.Value
# setter auto-coerces such strings to numbers, so for every cell we briefly
# force a text NumberFormat before the write, then ClearFormats() to drop the
# now-unneeded explicit style (keeping the cell's original default styling)
# while leaving the stored value as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue 'D2' '34.334.73'
Set-TextValue 'E2' '  +0.47%  '
Set-TextValue 'D3' '1.790.82'
Set-TextValue 'E3' '  +0.16%  '
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '226.49'
Set-TextValue 'E5' '  -0.04%  '
Set-TextValue 'E6' '  +1.63%  '
Set-TextValue 'E7' '  -0.01%  '
Set-TextValue 'D8' '32.24'
Set-TextValue 'E8' '  +0.61%  '
Set-TextValue 'E9' '  +0.88%  '
Set-TextValue 'D10' '0.0689'
Set-TextValue 'E10' '  -0.42%  '
Set-TextValue 'E11' '  +0.49%  '
Set-TextValue 'D12' '2.050.02'
Set-TextValue 'E12' '  +0.14%  '
Set-TextValue 'B13' 'Chainlink'
Set-TextValue 'C13' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D13' '10.94'
Set-TextValue 'E13' '  -3.56%  '
Set-TextValue 'B14' 'WrappedEther'
Set-TextValue 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D14' '1.776.16'
Set-TextValue 'E14' '  -0.87%  '
Set-TextValue 'E15' '  +0.87%  '
Set-TextValue 'D16' '34.313.48'
Set-TextValue 'E16' '  +0.55%  '
Set-TextValue 'E17' '  -0.06%  '
Set-TextValue 'D18' '68.07'
Set-TextValue 'E18' '  +0.09%  '
Set-TextValue 'D19' '0.0₃0801'
Set-TextValue 'E19' '  +2.52%  '
Set-TextValue 'D20' '246.78'
Set-TextValue 'E20' '  +0.43%  '
Set-TextValue 'D21' '10.93'
Set-TextValue 'E21' '  +0.67%  '
Set-TextValue 'E22' '  -0.09%  '
Set-TextValue 'E23' '  +1.23%  '
Set-TextValue 'D24' '2.05'
Set-TextValue 'E24' '  -0.16%  '
Set-TextValue 'D25' '162.32'
Set-TextValue 'E25' '  +0.45%  '
Set-TextValue 'E26' '  +0.23%  '
Set-TextValue 'D27' '16.36'
Set-TextValue 'E27' '  +0.18%  '
Set-TextValue 'E28' '  +1.37%  '
Set-TextValue 'E29' '  +0.06%  '
Set-TextValue 'E30' '  -0.49%  '
Set-TextValue 'B31' 'Hedera'
Set-TextValue 'C31' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D31' '0.0520'
Set-TextValue 'E31' '  +0.09%  '
Set-TextValue 'B32' 'Filecoin'
Set-TextValue 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D32' '3.78'
Set-TextValue 'E32' '  +2.96%  '
Set-TextValue 'B33' 'InternetComputer(DFINITY)'
Set-TextValue 'C33' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D33' '3.89'
Set-TextValue 'E33' '  +7.68%  '
Set-TextValue 'E34' '  -1.06%  '
Set-TextValue 'D35' '1.438.89'
Set-TextValue 'E35' '  -0.36%  '
Set-TextValue 'E36' '  +7.04%  '
Set-TextValue 'D37' '0.659'
Set-TextValue 'E37' '  +1.84%  '
Set-TextValue 'D38' '1.05'
Set-TextValue 'E38' '  +1.57%  '
Set-TextValue 'E39' '  -0.71%  '
Set-TextValue 'D40' '82.16'
Set-TextValue 'E40' '  +2.44%  '
Set-TextValue 'E41' '  +1.54%  '
Set-TextValue 'D42' '14.03'
Set-TextValue 'E42' '  +4.73%  '
Set-TextValue 'E43' '  +1.89%  '
Set-TextValue 'D44' '0.927'
Set-TextValue 'E44' '  +0.59%  '
Set-TextValue 'D45' '0.0521'
Set-TextValue 'E45' '  +1.92%  '
Set-TextValue 'D46' '6.04'
Set-TextValue 'E46' '  -0.39%  '
Set-TextValue 'E47' '  +0.23%  '
Set-TextValue 'D48' '1.945.77'
Set-TextValue 'E48' '  -0.13%  '
Set-TextValue 'D49' '105.31'
Set-TextValue 'E49' '  -2.39%  '
Set-TextValue 'D50' '0.0₆0131'
Set-TextValue 'E50' '  -5.99%  '
Set-TextValue 'E51' '  +0.00%  '
